$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 658.41174
$ws.Range("J28").Value = 1090
$ws.Range("L28").Value = 1090
$ws.Range("N28").Value = -2060
$ws.Range("H80").Value = 837.875
$ws.Range("I80").Value = 599.3333
$ws.Range("K80").Value = 1797.9999
$ws.Range("M80").Value = -799.9999
$ws.Range("H82").Value = 18108.25
$ws.Range("I82").Value = 795.5
$ws.Range("K82").Value = 2386.5
$ws.Range("M82").Value = -1980.5
$ws.Range("H83").Value = 837.875
$ws.Range("I83").Value = 599.3333
$ws.Range("K83").Value = 5393.9997
$ws.Range("M83").Value = -401.9997000000003
$ws.Range("H85").Value = 18108.25
$ws.Range("I85").Value = 795.5
$ws.Range("K85").Value = 2386.5
$ws.Range("M85").Value = -982.5
$ws.Range("H86").Value = 2879805.2
$ws.Range("J86").Value = 15600
$ws.Range("L86").Value = 15600
$ws.Range("N86").Value = -17846
$ws.Range("H87").Value = 129000
$ws.Range("J87").Value = 129000
$ws.Range("L87").Value = 129000
$ws.Range("N87").Value = -131496
$ws.Range("H88").Value = 742942.25
$ws.Range("I88").Value = 2459
$ws.Range("K88").Value = 2459
$ws.Range("M88").Value = -2053
$ws.Range("H89").Value = 2879805.2
$ws.Range("J89").Value = 15600
$ws.Range("L89").Value = 78000
$ws.Range("N89").Value = -89232
$ws.Range("H90").Value = 129000
$ws.Range("J90").Value = 129000
$ws.Range("L90").Value = 387000
$ws.Range("N90").Value = -399480
$ws.Range("H91").Value = 742942.25
$ws.Range("I91").Value = 2459
$ws.Range("K91").Value = 2459
$ws.Range("M91").Value = -1055
$ws.Range("H93").Value = 120000
$ws.Range("J93").Value = 120000
$ws.Range("L93").Value = 120000
$ws.Range("N93").Value = -124992
$ws.Range("H96").Value = 2169.842
$ws.Range("I96").Value = 1972.5714
$ws.Range("K96").Value = 5917.7142
$ws.Range("M96").Value = -4544.7142
$ws.Range("H98").Value = 76962130
$ws.Range("I98").Value = 83375610
$ws.Range("K98").Value = 83375610
$ws.Range("M98").Value = -83374112
$ws.Range("H105").Value = 57497.5
$ws.Range("J105").Value = 57497.5
$ws.Range("L105").Value = 57497.5
$ws.Range("N105").Value = -64485.5
$ws.Range("H112").Value = 1161.5385
$ws.Range("I112").Value = 460
$ws.Range("J112").Value = 3500
$ws.Range("K112").Value = 1380
$ws.Range("L112").Value = 10500
$ws.Range("M112").Value = -272
$ws.Range("N112").Value = -12716
$ws.Range("H116").Value = 6415.7144
$ws.Range("I116").Value = 5985
$ws.Range("K116").Value = 5985
$ws.Range("M116").Value = -2543
$ws.Range("H122").Value = 76962130
$ws.Range("I122").Value = 83375610
$ws.Range("K122").Value = 250126830
$ws.Range("M122").Value = -250124380
$ws.Range("H137").Value = 3079
$ws.Range("I137").Value = 1695
$ws.Range("J137").Value = 9999
$ws.Range("K137").Value = 5085
$ws.Range("L137").Value = 29997
$ws.Range("M137").Value = -2535
$ws.Range("N137").Value = -35097
$ws.Range("H138").Value = 2933.3218
$ws.Range("I138").Value = 1509.579
$ws.Range("J138").Value = 3331.1323
$ws.Range("K138").Value = 4528.737
$ws.Range("L138").Value = 9993.3969
$ws.Range("M138").Value = 611.2629999999999
$ws.Range("N138").Value = -20273.3969
$ws.Range("H141").Value = 2036.4584
$ws.Range("J141").Value = 2844.6667
$ws.Range("L141").Value = 8534.000100000001
$ws.Range("N141").Value = -18894.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9261583
$ws.Range("I32").Value = 10205482
$ws.Range("K32").Value = 10205482
$ws.Range("M32").Value = -10205195
$ws.Range("H61").Value = 83514080
$ws.Range("I61").Value = 166680000
$ws.Range("K61").Value = 166680000
$ws.Range("M61").Value = -166679788
$ws.Range("H105").Value = 70370
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H112").Value = 83878
$ws.Range("J112").Value = 83878
$ws.Range("L112").Value = 83878
$ws.Range("N112").Value = -86832
$ws.Range("H136").Value = 83514080
$ws.Range("I136").Value = 166680000
$ws.Range("K136").Value = 500040000
$ws.Range("M136").Value = -500037450

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5019.25
$ws.Range("I20").Value = 5157.364
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 5157.364
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -4910.364
$ws.Range("N20").Value = -3994
$ws.Range("H38").Value = 38501
$ws.Range("J38").Value = 38501
$ws.Range("L38").Value = 38501
$ws.Range("N38").Value = -39333
$ws.Range("H106").Value = 45789.6
$ws.Range("J106").Value = 45789.6
$ws.Range("L106").Value = 45789.6
$ws.Range("N106").Value = -48313.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1980
$ws.Range("I58").Value = 1980
$ws.Range("K58").Value = 1980
$ws.Range("M58").Value = -1777
$ws.Range("H64").Value = 49995
$ws.Range("J64").Value = 49995
$ws.Range("L64").Value = 49995
$ws.Range("N64").Value = -50491
$ws.Range("H67").Value = 49995
$ws.Range("J67").Value = 49995
$ws.Range("L67").Value = 49995
$ws.Range("N67").Value = -51711
$ws.Range("H136").Value = 1980
$ws.Range("I136").Value = 1980
$ws.Range("K136").Value = 5940
$ws.Range("M136").Value = -3390

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 95.76922999999999
$ws.Range("I2").Value = 59.2
$ws.Range("K2").Value = 355.2
$ws.Range("M2").Value = -242.2
$ws.Range("H86").Value = 588.375
$ws.Range("J86").Value = 611.4
$ws.Range("L86").Value = 1834.2
$ws.Range("N86").Value = -4206.2
$ws.Range("H87").Value = 37507
$ws.Range("I87").Value = 37507
$ws.Range("K87").Value = 112521
$ws.Range("M87").Value = -111273
$ws.Range("H88").Value = 3866.6667
$ws.Range("I88").Value = 2800
$ws.Range("K88").Value = 8400
$ws.Range("M88").Value = -7972
$ws.Range("H89").Value = 588.375
$ws.Range("J89").Value = 611.4
$ws.Range("L89").Value = 5502.599999999999
$ws.Range("N89").Value = -17358.6
$ws.Range("H90").Value = 37507
$ws.Range("I90").Value = 37507
$ws.Range("K90").Value = 337563
$ws.Range("M90").Value = -331323
$ws.Range("H91").Value = 3866.6667
$ws.Range("I91").Value = 2800
$ws.Range("K91").Value = 8400
$ws.Range("M91").Value = -6918
$ws.Range("H107").Value = 547.43475
$ws.Range("J107").Value = 618.3570999999999
$ws.Range("L107").Value = 1855.0713
$ws.Range("N107").Value = -5695.0713
$ws.Range("H139").Value = 2503.5557
$ws.Range("J139").Value = 3666.6667
$ws.Range("L139").Value = 11000.0001
$ws.Range("N139").Value = -21280.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5731.5
$ws.Range("I70").Value = 5731.5
$ws.Range("K70").Value = 5731.5
$ws.Range("M70").Value = -5461.5
$ws.Range("H73").Value = 5731.5
$ws.Range("I73").Value = 5731.5
$ws.Range("K73").Value = 5731.5
$ws.Range("M73").Value = -4795.5
$ws.Range("H80").Value = 3886.875
$ws.Range("I80").Value = 3619.2
$ws.Range("J80").Value = 4333
$ws.Range("K80").Value = 3619.2
$ws.Range("L80").Value = 4333
$ws.Range("M80").Value = -2621.2
$ws.Range("N80").Value = -6329
$ws.Range("H83").Value = 3886.875
$ws.Range("I83").Value = 3619.2
$ws.Range("J83").Value = 4333
$ws.Range("K83").Value = 18096
$ws.Range("L83").Value = 21665
$ws.Range("M83").Value = -13104
$ws.Range("N83").Value = -31649
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 117450.5
$ws.Range("J110").Value = 117450.5
$ws.Range("L110").Value = 117450.5
$ws.Range("N110").Value = -125630.5
$ws.Range("H132").Value = 202850.9
$ws.Range("J132").Value = 252751.25
$ws.Range("L132").Value = 758253.75
$ws.Range("N132").Value = -763313.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3037.5
$ws.Range("I132").Value = 2516.6667
$ws.Range("K132").Value = 7550.000100000001
$ws.Range("M132").Value = -5020.000100000001
$ws.Range("H136").Value = 1311
$ws.Range("J136").Value = 1502.5
$ws.Range("L136").Value = 4507.5
$ws.Range("N136").Value = -9607.5
